$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.292
$ws.Range("A6").Value = -22.291
$ws.Range("A7").Value = -20.053
$ws.Range("E7").Value = 16.327
$ws.Range("E12").Value = 17.638
$ws.Range("E15").Value = 16.113
$ws.Range("A16").Value = -21.923
$ws.Range("A20").Value = -20.31
$ws.Range("E20").Value = 16.197
$ws.Range("E21").Value = 16.534
$ws.Range("E22").Value = 16.587
$ws.Range("E23").Value = 16.434
$ws.Range("A28").Value = -22.016
$ws.Range("A29").Value = -21.344
$ws.Range("E29").Value = 16.97
$ws.Range("A32").Value = -21.729
$ws.Range("E34").Value = 16.747
$ws.Range("A40").Value = -19.896
$ws.Range("E42").Value = 16.539
$ws.Range("E43").Value = 16.917
$ws.Range("E44").Value = 16.535
$ws.Range("E45").Value = 16.774
$ws.Range("A46").Value = -21.816
$ws.Range("E46").Value = 16.679
$ws.Range("E50").Value = 16.403
$ws.Range("A51").Value = -21.934
$ws.Range("E51").Value = 16.881
$ws.Range("A52").Value = -22.084
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.395
$ws.Range("A62").Value = -22.085
$ws.Range("A66").Value = -21.546
$ws.Range("E66").Value = 17.43
$ws.Range("E67").Value = 17.571
$ws.Range("A73").Value = -20.213
$ws.Range("A74").Value = -21.112
$ws.Range("E79").Value = 16.957
$ws.Range("E84").Value = 16.606
$ws.Range("A92").Value = -21.515
$ws.Range("E92").Value = 17.83
$ws.Range("E97").Value = 16.822
$ws.Range("A100").Value = -22.217
